# fix: remove cohorts from tablename and remove duplicate recruit_age
#
# The "Variables" sheet had "recruit_age" listed twice: once at row 105
# (the original/duplicate entry) and again near the bottom of the sheet
# (row 314, alongside the other recently-appended variables such as
# cats_preg/cats_quant_preg/dogs_preg/dogs_quant_preg). This removes the
# duplicate row 105 entry; Excel shifts every row below it up by one,
# which is exactly what the published diff shows (row 106 "sleept_psc"
# becomes the new row 105, etc., all the way down to the tail rows).

$wb = $excel.ActiveWorkbook
$wsVariables = $wb.Worksheets.Item("Variables")
$wsCategories = $wb.Worksheets.Item("Categories")

# Delete the duplicate "recruit_age" row (row 105) on the Variables sheet.
$wsVariables.Rows.Item(105).Delete()

# Reflect the end state of the edit: the workbook was left with the
# Variables sheet active/selected, with the (now empty, former row-105)
# row band A105:XFD105 highlighted -- consistent with how Excel leaves
# the selection right after a "Delete Row" operation.
$wsVariables.Activate() | Out-Null
$wsVariables.Range("A105:XFD105").Select() | Out-Null
